$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 898.2727
$ws.Cells.Item(40, 9).Value = 880.1111
$ws.Cells.Item(40, 10).Value = 980
$ws.Cells.Item(40, 11).Value = 880.1111
$ws.Cells.Item(40, 12).Value = 980
$ws.Cells.Item(40, 13).Value = -705.1111
$ws.Cells.Item(40, 14).Value = -1330
$ws.Cells.Item(80, 8).Value = 406.25
$ws.Cells.Item(80, 9).Value = 220.76923
$ws.Cells.Item(80, 11).Value = 662.30769
$ws.Cells.Item(80, 13).Value = 335.69231
$ws.Cells.Item(83, 8).Value = 406.25
$ws.Cells.Item(83, 9).Value = 220.76923
$ws.Cells.Item(83, 11).Value = 1986.92307
$ws.Cells.Item(83, 13).Value = 3005.07693
$ws.Cells.Item(96, 8).Value = 1793.091
$ws.Cells.Item(96, 9).Value = 2330.7144
$ws.Cells.Item(96, 10).Value = 852.25
$ws.Cells.Item(96, 11).Value = 6992.1432
$ws.Cells.Item(96, 12).Value = 2556.75
$ws.Cells.Item(96, 13).Value = -5619.1432
$ws.Cells.Item(96, 14).Value = -5302.75
$ws.Cells.Item(112, 8).Value = 2005.2778
$ws.Cells.Item(112, 10).Value = 2193.1333
$ws.Cells.Item(112, 12).Value = 6579.3999
$ws.Cells.Item(112, 14).Value = -8795.3999
$ws.Cells.Item(137, 8).Value = 2732.7273
$ws.Cells.Item(137, 9).Value = 2030.75
$ws.Cells.Item(137, 11).Value = 6092.25
$ws.Cells.Item(137, 13).Value = -3542.25
$ws.Cells.Item(138, 8).Value = 2559.694
$ws.Cells.Item(138, 9).Value = 1530.4375
$ws.Cells.Item(138, 10).Value = 2798.3623
$ws.Cells.Item(138, 11).Value = 4591.3125
$ws.Cells.Item(138, 12).Value = 8395.086899999998
$ws.Cells.Item(138, 13).Value = 548.6875
$ws.Cells.Item(138, 14).Value = -18675.0869
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 111112664
$ws.Cells.Item(61, 10).Value = 3000
$ws.Cells.Item(61, 12).Value = 3000
$ws.Cells.Item(61, 14).Value = -3424
$ws.Cells.Item(74, 8).Value = 1998.5862
$ws.Cells.Item(74, 9).Value = 1107.7894
$ws.Cells.Item(74, 11).Value = 1107.7894
$ws.Cells.Item(74, 13).Value = -233.7893999999999
$ws.Cells.Item(77, 8).Value = 1998.5862
$ws.Cells.Item(77, 9).Value = 1107.7894
$ws.Cells.Item(77, 11).Value = 5538.946999999999
$ws.Cells.Item(77, 13).Value = -1170.946999999999
$ws.Cells.Item(97, 8).Value = 17321.666
$ws.Cells.Item(97, 9).Value = 786
$ws.Cells.Item(97, 10).Value = 100000
$ws.Cells.Item(97, 11).Value = 786
$ws.Cells.Item(97, 12).Value = 100000
$ws.Cells.Item(97, 13).Value = -290
$ws.Cells.Item(97, 14).Value = -100992
$ws.Cells.Item(136, 8).Value = 111112664
$ws.Cells.Item(136, 10).Value = 3000
$ws.Cells.Item(136, 12).Value = 9000
$ws.Cells.Item(136, 14).Value = -14100
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 457.6
$ws.Cells.Item(22, 9).Value = 350
$ws.Cells.Item(22, 10).Value = 529.3333
$ws.Cells.Item(22, 11).Value = 350
$ws.Cells.Item(22, 12).Value = 529.3333
$ws.Cells.Item(22, 13).Value = -177
$ws.Cells.Item(22, 14).Value = -875.3333
$ws.Cells.Item(134, 8).Value = 3108.2954
$ws.Cells.Item(134, 9).Value = 694.2727
$ws.Cells.Item(134, 11).Value = 2082.8181
$ws.Cells.Item(134, 13).Value = 452.1819
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 211.72728
$ws.Cells.Item(7, 9).Value = 103.71429
$ws.Cells.Item(7, 10).Value = 400.75
$ws.Cells.Item(7, 11).Value = 103.71429
$ws.Cells.Item(7, 12).Value = 400.75
$ws.Cells.Item(7, 13).Value = 9.285709999999995
$ws.Cells.Item(7, 14).Value = -626.75
$ws.Cells.Item(19, 8).Value = 112.25
$ws.Cells.Item(19, 9).Value = 112.25
$ws.Cells.Item(19, 11).Value = 112.25
$ws.Cells.Item(19, 13).Value = 57.75
$ws.Cells.Item(24, 8).Value = 112.25
$ws.Cells.Item(24, 9).Value = 112.25
$ws.Cells.Item(24, 11).Value = 112.25
$ws.Cells.Item(24, 13).Value = 57.75
$ws.Cells.Item(99, 8).Value = 1549285.2
$ws.Cells.Item(99, 9).Value = 2632745
$ws.Cells.Item(99, 11).Value = 2632745
$ws.Cells.Item(99, 13).Value = -2631247
$ws.Cells.Item(107, 8).Value = 1016.45
$ws.Cells.Item(107, 9).Value = 654.2857
$ws.Cells.Item(107, 10).Value = 1861.5
$ws.Cells.Item(107, 11).Value = 654.2857
$ws.Cells.Item(107, 12).Value = 1861.5
$ws.Cells.Item(107, 13).Value = 1265.7143
$ws.Cells.Item(107, 14).Value = -5701.5
$ws.Cells.Item(126, 8).Value = 1549285.2
$ws.Cells.Item(126, 9).Value = 2632745
$ws.Cells.Item(126, 11).Value = 7898235
$ws.Cells.Item(126, 13).Value = -7895765
$ws.Cells.Item(132, 8).Value = 1989.4166
$ws.Cells.Item(132, 9).Value = 1619.7059
$ws.Cells.Item(132, 10).Value = 2887.2856
$ws.Cells.Item(132, 11).Value = 4859.1177
$ws.Cells.Item(132, 12).Value = 8661.856800000001
$ws.Cells.Item(132, 13).Value = -2329.1177
$ws.Cells.Item(132, 14).Value = -13721.8568
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 481.45
$ws.Cells.Item(5, 9).Value = 374.1111
$ws.Cells.Item(5, 10).Value = 1447.5
$ws.Cells.Item(5, 11).Value = 1122.3333
$ws.Cells.Item(5, 12).Value = 4342.5
$ws.Cells.Item(5, 13).Value = -1010.3333
$ws.Cells.Item(5, 14).Value = -4566.5
$ws.Cells.Item(113, 8).Value = 709.7659
$ws.Cells.Item(113, 9).Value = 649.17645
$ws.Cells.Item(113, 10).Value = 744.1
$ws.Cells.Item(113, 11).Value = 1947.52935
$ws.Cells.Item(113, 12).Value = 2232.3
$ws.Cells.Item(113, 13).Value = 222.4706499999998
$ws.Cells.Item(113, 14).Value = -6572.3
$ws.Cells.Item(135, 8).Value = 481.45
$ws.Cells.Item(135, 9).Value = 374.1111
$ws.Cells.Item(135, 10).Value = 1447.5
$ws.Cells.Item(135, 11).Value = 3366.9999
$ws.Cells.Item(135, 12).Value = 13027.5
$ws.Cells.Item(135, 13).Value = -831.9999000000003
$ws.Cells.Item(135, 14).Value = -18097.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1214.9286
$ws.Cells.Item(97, 9).Value = 1140.8
$ws.Cells.Item(97, 10).Value = 1400.25
$ws.Cells.Item(97, 11).Value = 1140.8
$ws.Cells.Item(97, 12).Value = 1400.25
$ws.Cells.Item(97, 13).Value = -644.8
$ws.Cells.Item(97, 14).Value = -2392.25
$ws.Cells.Item(98, 8).Value = 12643
$ws.Cells.Item(98, 10).Value = 12643
$ws.Cells.Item(98, 12).Value = 12643
$ws.Cells.Item(98, 14).Value = -18633
$ws.Cells.Item(107, 8).Value = 3205662
$ws.Cells.Item(107, 9).Value = 4808318
$ws.Cells.Item(107, 11).Value = 4808318
$ws.Cells.Item(107, 13).Value = -4806398
$ws.Cells.Item(126, 8).Value = 1817
$ws.Cells.Item(126, 9).Value = 1594.3636
$ws.Cells.Item(126, 11).Value = 4783.0908
$ws.Cells.Item(126, 13).Value = -2313.0908
$ws.Cells.Item(132, 8).Value = 7216.9585
$ws.Cells.Item(132, 9).Value = 9592.5
$ws.Cells.Item(132, 10).Value = 3891.2
$ws.Cells.Item(132, 11).Value = 28777.5
$ws.Cells.Item(132, 12).Value = 11673.6
$ws.Cells.Item(132, 13).Value = -26247.5
$ws.Cells.Item(132, 14).Value = -16733.6
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1555.5
$ws.Cells.Item(22, 9).Value = 1111
$ws.Cells.Item(22, 10).Value = 2000
$ws.Cells.Item(22, 11).Value = 1111
$ws.Cells.Item(22, 12).Value = 2000
$ws.Cells.Item(22, 13).Value = -816
$ws.Cells.Item(22, 14).Value = -2590
$ws.Cells.Item(27, 8).Value = 1555.5
$ws.Cells.Item(27, 9).Value = 1111
$ws.Cells.Item(27, 10).Value = 2000
$ws.Cells.Item(27, 11).Value = 1111
$ws.Cells.Item(27, 12).Value = 2000
$ws.Cells.Item(27, 13).Value = -1004
$ws.Cells.Item(27, 14).Value = -2214
$ws.Cells.Item(68, 8).Value = 1286.5
$ws.Cells.Item(68, 9).Value = 1286.5
$ws.Cells.Item(68, 11).Value = 1286.5
$ws.Cells.Item(68, 13).Value = -537.5
$ws.Cells.Item(71, 8).Value = 1286.5
$ws.Cells.Item(71, 9).Value = 1286.5
$ws.Cells.Item(71, 11).Value = 6432.5
$ws.Cells.Item(71, 13).Value = -2688.5
$ws.Cells.Item(132, 8).Value = 2210.3142
$ws.Cells.Item(132, 9).Value = 1778.52
$ws.Cells.Item(132, 11).Value = 5335.559999999999
$ws.Cells.Item(132, 13).Value = -2805.559999999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 3750
$ws.Cells.Item(96, 9).Value = 6000
$ws.Cells.Item(96, 10).Value = 1500
$ws.Cells.Item(96, 11).Value = 6000
$ws.Cells.Item(96, 12).Value = 1500
$ws.Cells.Item(96, 13).Value = -4627
$ws.Cells.Item(96, 14).Value = -4246
$ws.Cells.Item(132, 8).Value = 2578.5095
$ws.Cells.Item(132, 9).Value = 2587.9768
$ws.Cells.Item(132, 11).Value = 7763.930399999999
$ws.Cells.Item(132, 13).Value = -5233.930399999999

Write-Output "Applied 192 cell updates across 8 sheets"